# Update column G ("K") values for rows 2-28 on the active sheet.
# These are the recalculated strikeout/K counts (s_vals) for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 7
    4  = 1
    5  = 3
    6  = 7
    7  = 3
    8  = 6
    9  = 2
    10 = 3
    11 = 0
    12 = 6
    13 = 5
    14 = 8
    15 = 5
    16 = 9
    17 = 4
    18 = 2
    19 = 5
    20 = 2
    21 = 1
    22 = 7
    23 = 7
    24 = 5
    25 = 4
    26 = 3
    27 = 3
    28 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
